# Applies the cryptos.xlsx data refresh described by the diff
# (GitHub Actions scheduled crypto price/volume update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.909.01"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.667.35"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'215.46"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'0.534"
$ws.Range("E6").Value = "  +5.02%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'20.24"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Value = "'0.0895"
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "1.902.26"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "1.659.75"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'0.525"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'66.12"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "26.925.88"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "'234.76"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'7.97"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.15"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  -2.10%  "
$ws.Range("D25").Value = "'146.18"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "'7.13"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "1.454.25"
$ws.Range("E33").Value = "  -3.99%  "
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("E35").Value = "  +2.06%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("D38").Value = "'0.903"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").Value = "'5.71"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D43").Value = "'66.14"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'0.969"
$ws.Range("E44").Value = "  +5.36%  "
$ws.Range("D45").Value = "1.810.21"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'0.783"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "'90.54"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("E50").Value = "  +4.31%  "
$ws.Range("E51").Value = "  +0.01%  "
